# Weekly update: three new price records (Damasco, Castle Brite, week of
# 2021-12-09 / serial 44539, "Provincia de San Felipe de Aconcagua") are
# inserted at the top of the data block (row 14), pushing the existing
# rows 14-22 down to rows 17-25 unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at row 14 (shifts old rows 14-22 down to 17-25)
$ws.Range("A14:A16").EntireRow.Insert()

# Row 14 - Damasco, Castle Brite, Especial
$ws.Range("A14").Value = 8
$ws.Range("B14").Value = "Terminal La Palmera de La Serena"
$ws.Range("C14").Value = "Coquimbo"
$ws.Range("D14").Value = 44539
$ws.Range("E14").Value = 4
$ws.Range("F14").Value = "Fruta"
$ws.Range("G14").Value = 100103
$ws.Range("H14").Value = "Frutos de hueso (carozo)"
$ws.Range("I14").Value = 100103003
$ws.Range("J14").Value = "Damasco"
$ws.Range("K14").Value = "Castle Brite"
$ws.Range("L14").Value = "Especial"
$ws.Range("M14").Value = 160
$ws.Range("N14").Value = 24500
$ws.Range("O14").Value = 25000
$ws.Range("P14").Value = 24750
$ws.Range("Q14").Value = "`$/caja 15 kilos"
$ws.Range("R14").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S14").Value = 1650
$ws.Range("T14").Value = 15

# Row 15 - Damasco, Castle Brite, Primera
$ws.Range("A15").Value = 8
$ws.Range("B15").Value = "Terminal La Palmera de La Serena"
$ws.Range("C15").Value = "Coquimbo"
$ws.Range("D15").Value = 44539
$ws.Range("E15").Value = 4
$ws.Range("F15").Value = "Fruta"
$ws.Range("G15").Value = 100103
$ws.Range("H15").Value = "Frutos de hueso (carozo)"
$ws.Range("I15").Value = 100103003
$ws.Range("J15").Value = "Damasco"
$ws.Range("K15").Value = "Castle Brite"
$ws.Range("L15").Value = "Primera"
$ws.Range("M15").Value = 160
$ws.Range("N15").Value = 22500
$ws.Range("O15").Value = 23000
$ws.Range("P15").Value = 22750
$ws.Range("Q15").Value = "`$/caja 15 kilos"
$ws.Range("R15").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S15").Value = 1517
$ws.Range("T15").Value = 15

# Row 16 - Damasco, Castle Brite, Segunda
$ws.Range("A16").Value = 8
$ws.Range("B16").Value = "Terminal La Palmera de La Serena"
$ws.Range("C16").Value = "Coquimbo"
$ws.Range("D16").Value = 44539
$ws.Range("E16").Value = 4
$ws.Range("F16").Value = "Fruta"
$ws.Range("G16").Value = 100103
$ws.Range("H16").Value = "Frutos de hueso (carozo)"
$ws.Range("I16").Value = 100103003
$ws.Range("J16").Value = "Damasco"
$ws.Range("K16").Value = "Castle Brite"
$ws.Range("L16").Value = "Segunda"
$ws.Range("M16").Value = 200
$ws.Range("N16").Value = 18000
$ws.Range("O16").Value = 18500
$ws.Range("P16").Value = 18250
$ws.Range("Q16").Value = "`$/caja 15 kilos"
$ws.Range("R16").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S16").Value = 1217
$ws.Range("T16").Value = 15
